$d = $word.ActiveDocument

# 1. Replace "studying difficulties" with "some personal reasons" in the
#    introductory paragraph.
$d.Content.Find.Execute("studying difficulties", $true, $false, $false, `
    $false, $false, $true, 1, $false, "some personal reasons", 2) | Out-Null

# 2. Locate the point right after the just-inserted replacement text; this is
#    where the new "_GoBack" bookmark needs to live (splitting the run in two,
#    matching Word's own behaviour when the cursor is left at that spot).
$splitPoint = $d.Content
$splitPoint.Find.Execute("some personal reasons", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitPoint.Collapse(0)   # wdCollapseEnd

# 3. Remove the existing "_GoBack" bookmark (currently sitting by itself in an
#    empty paragraph further down the document).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 4. Re-create the "_GoBack" bookmark at the new split point.
$d.Bookmarks.Add("_GoBack", $splitPoint) | Out-Null
